# Insert a new weekly record for "Macroferia Regional de Talca" (Piña,
# Caramelo, Segunda) dated 2021-11-04 (serial 44504) right above the
# existing 2021-10-18 (44487) record, pushing the rest of the table down
# by one row (old row 140 -> new row 141, ..., old row 167 -> new row 168).

$wb2 = $excel.ActiveWorkbook
$ws  = $wb2.ActiveSheet

# Shift rows 140:167 down to 141:168, leaving a blank row 140 to fill in.
$ws.Rows.Item(140).Insert()

$newRow = 140

$ws.Cells.Item($newRow, 1).Value  = 5
$ws.Cells.Item($newRow, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value  = "Maule"
$ws.Cells.Item($newRow, 4).Value  = 44504
$ws.Cells.Item($newRow, 5).Value  = 7
$ws.Cells.Item($newRow, 6).Value  = "Fruta"
$ws.Cells.Item($newRow, 7).Value  = 100108
$ws.Cells.Item($newRow, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value  = 100108005
$ws.Cells.Item($newRow, 10).Value = "Piña"
$ws.Cells.Item($newRow, 11).Value = "Caramelo"
$ws.Cells.Item($newRow, 12).Value = "Segunda"
$ws.Cells.Item($newRow, 13).Value = 210
$ws.Cells.Item($newRow, 14).Value = 18000
$ws.Cells.Item($newRow, 15).Value = 18000
$ws.Cells.Item($newRow, 16).Value = 18000
$ws.Cells.Item($newRow, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item($newRow, 18).Value = "Ecuador"
$ws.Cells.Item($newRow, 19).Value = 1286
$ws.Cells.Item($newRow, 20).Value = 14
